# Excel import template update:
#  - the "Surname"/patronymic data column (E) now stores the student's
#    gender (М/Ж) instead of a patronymic, for rows 2-4
#  - selection cursor left on H9 (last place the editor clicked)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "М"
$ws.Range("E3").Value = "Ж"
$ws.Range("E4").Value = "М"

$null = $ws.Range("H9").Select()
